$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.05"
$ws.Range("E2").Value = "'1.43%"
$ws.Range("D3").Value = "'32.14"
$ws.Range("E3").Value = "'2.90%"
$ws.Range("D4").Value = "'4.932"
$ws.Range("E4").Value = "'-3.60%"
$ws.Range("D5").Value = "'0.07828"
$ws.Range("E5").Value = "'-1.60%"
$ws.Range("D6").Value = "'2.040"
$ws.Range("E6").Value = "'-8.49%"
$ws.Range("D7").Value = "'7.836"
$ws.Range("E7").Value = "'0.80%"
$ws.Range("D8").Value = "'3.823"
$ws.Range("E8").Value = "'-1.04%"
$ws.Range("D9").Value = "'0.9224"
$ws.Range("E9").Value = "'-0.15%"
$ws.Range("D10").Value = "'0.1762"
$ws.Range("E10").Value = "'1.95%"
$ws.Range("D11").Value = "'0.07902"
$ws.Range("E11").Value = "'6.08%"
$ws.Range("D12").Value = "'0.08599"
$ws.Range("E12").Value = "'-8.17%"
$ws.Range("D13").Value = "'0.03160"
$ws.Range("E13").Value = "'3.18%"
$ws.Range("D14").Value = "'0.1005"
$ws.Range("E14").Value = "'0.06%"
$ws.Range("D15").Value = "'0.001513"
$ws.Range("E15").Value = "'0.05%"
$ws.Range("D16").Value = "'0.005880"
$ws.Range("E16").Value = "'-2.05%"
$ws.Range("E17").Value = "'2,110.40%"
$ws.Range("E18").Value = "'-0.42%"
$ws.Range("D19").Value = "'2.157"
$ws.Range("E19").Value = "'-4.90%"
$ws.Range("D20").Value = "'0.3309"
$ws.Range("E20").Value = "'1.15%"
$ws.Range("D21").Value = "'0.1318"
$ws.Range("E21").Value = "'-1.18%"
$ws.Range("D22").Value = "'4.270"
$ws.Range("E22").Value = "'9.23%"
$ws.Range("D23").Value = "'0.1991"
$ws.Range("E23").Value = "'17.23%"
$ws.Range("D24").Value = "'0.04574"
$ws.Range("E24").Value = "'-1.03%"
$ws.Range("E25").Value = "'-1.84%"
$ws.Range("D26").Value = "'0.004446"
$ws.Range("E26").Value = "'-0.73%"
$ws.Range("E27").Value = "'4.31%"
$ws.Range("D39").Value = "'0.01741"
$ws.Range("E39").Value = "'-1.04%"
$ws.Range("D40").Value = "'0.04779"
$ws.Range("E40").Value = "'3.72%"
$ws.Range("D41").Value = "'0.007560"
$ws.Range("E41").Value = "'8.45%"
$ws.Range("E42").Value = "'0.23%"
$ws.Range("D43").Value = "'0.002361"
$ws.Range("E43").Value = "'7.91%"
$ws.Range("E44").Value = "'4.83%"
$ws.Range("D45").Value = "'0.00006302"
$ws.Range("E45").Value = "'-0.07%"
$ws.Range("E47").Value = "'-61.09%"
$ws.Range("D48").Value = "'0.8234"
$ws.Range("E48").Value = "'-28.88%"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D50").Value = "'0.0002001"
